{"js": "// Abstract paragraph wording update:\n//  - \"Snorkel surveys\" -> \"snorkel surveys\" (lowercase the stray capital S)\n//  - \"Crews survey 20 sampling sites\" -> \"Crews survey units within 20 sampling sections\"\nconst body = context.document.body;\n\nconst snorkelResults = body.search(\"Snorkel surveys\", { matchCase: true });\nsnorkelResults.load(\"text\");\n\nconst crewsResults = body.search(\"Crews survey 20 sampling sites\", { matchCase: true });\ncrewsResults.load(\"text\");\n\nawait context.sync();\n\nif (snorkelResults.items.length > 0) {\n  snorkelResults.items[0].insertText(\"snorkel surveys\", Word.InsertLocation.replace);\n}\n\nif (crewsResults.items.length > 0) {\n  crewsResults.items[0].insertText(\n    \"Crews survey units within 20 sampling sections\",\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n", "ps1": "# Abstract paragraph wording update:\n#  - \"Snorkel surveys\" -> \"snorkel surveys\" (lowercase the stray capital S)\n#  - \"Crews survey 20 sampling sites\" -> \"Crews survey units within 20 sampling sections\"\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"Snorkel surveys\"\n$find1.MatchCase = $true\n$find1.Replacement.Text = \"snorkel surveys\"\n$find1.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,$wdReplaceAll)\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Crews survey 20 sampling sites\"\n$find2.MatchCase = $true\n$find2.Replacement.Text = \"Crews survey units within 20 sampling sections\"\n$find2.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,$wdReplaceAll)\n"}
